$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bitacora")

# --- New row 8: "Creacion modelo conceptual iteracion 2" ---------------
# Column A holds the text "2.0" (same iteration label as row 7). Typing
# "2.0" directly gets auto-coerced to the number 2 (Excel's normal
# numeric-literal detection), so build it as a text formula in a scratch
# cell and paste-special the resulting value into A8 — that keeps it text
# without ever touching NumberFormat (which would otherwise leave a new,
# unused style entry behind in the workbook).
$scratch = $ws.Range("Z100")
$scratch.Formula = '="2.0"'
$scratch.Copy()
$ws.Range("A8").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("B8").Value = "H1 – Bicicletas"
$ws.Range("C8").Value = "feature/modelo_conceptual_iter_2"
$ws.Range("D8").Value = "conceptual"
$ws.Range("E8").Value = "andr4f"
$ws.Range("F8").Value = "andr4f, Angel Trillo, Yineth Avila"
$ws.Range("G8").Value = "https://github.com/andr4f/bici-go-bd/pull/13"
$ws.Range("H8").Value = "v0.2-iter2-conceptual"
$ws.Range("I8").Value = "Creacion modelo conceptual iteracion 2"

# Dates use the same "d-mmm" display format already used by J6:K7 (built
# in numFmtId 16), so this reuses the existing style instead of adding a
# new one. Assign the raw serial number first — typing a recognisable
# date STRING (e.g. "10/12/2025") makes Excel auto-apply its own default
# date format first, which would create a throwaway style entry before
# we get a chance to override it.
$ws.Range("J8").Value = 45942
$ws.Range("K8").Value = 45942
$ws.Range("J8:K8").NumberFormat = "d-mmm"

$ws.Range("L8").Value = "Completado"
$ws.Range("M8").Value = "bici_go_bd/modelos/conceptual/iteracion_2"

# Hyperlink for the PR URL in G8 (mirrors G4:G7).
$ws.Hyperlinks.Add($ws.Range("G8"), "https://github.com/andr4f/bici-go-bd/pull/13") | Out-Null
$ws.Range("G8").Style = "Hipervínculo"

# --- View state: drop the J1 scroll anchor, select G30:G31 -------------
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G30:G31").Select() | Out-Null
